{"js": "// Translate the French facilitator-guide labels to Italian.\n// Each (old, new) pair below corresponds to one <w:t> run in the document;\n// some French labels (\"Mat\u00e9riel\", \"Invitation \u00e0 la discussion\") occur twice\n// and both occurrences must be replaced.\nconst replacements = [\n  [\"Titre de la vid\u00e9o\", \"Titolo del video\"],\n  [\"Sujet\", \"Argomento\"],\n  [\"G\u00e9om\u00e9trie\", \"Geometria\"],\n  [\"Objectif(s)\", \"Obiettivo/i\"],\n  [\"Dur\u00e9e\", \"Lunghezza\"],\n  [\"Localisation du camp\", \"Posizione del Campo\"],\n  [\"Assistant\u00b7e\u00b7s\", \"Facilitatori\"],\n  [\"N. d'\u00e9tudiant\u00b7e\u00b7s\", \"N. di studenti\"],\n  [\"Date\", \"Data\"],\n  [\"Mat\u00e9riel\", \"Materiale\"],\n  [\"requis\", \"necessario\"],\n  [\"Pr\u00e9paration\", \"Preparazione\"],\n  [\"Temps vid\u00e9o\", \"Tempo del Video\"],\n  [\"Introduction de la vid\u00e9o\", \"Introduzione al video\"],\n  [\"Introduction de la premi\u00e8re activit\u00e9\", \"Introduzione all'attivit\u00e0 principale\"],\n  [\"Les 8 graphes dessin\u00e9s pour la vid\u00e9o\", \"Gli 8 grafici disegnati per il video\"],\n  [\"Introduction de l'exp\u00e9rience\", \"Introduzione all'esperimento\"],\n  [\"Invitation \u00e0 la discussion\", \"Invito alla discussione\"],\n  [\"Conclusion\", \"Conclusione\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    // Only replace runs whose full text equals the French label exactly,\n    // so we never touch a longer string that merely contains it.\n    if (range.text === oldText) {\n      range.insertText(newText, \"Replace\");\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the French facilitator-guide labels to Italian.\n# Each (old, new) pair corresponds to one <w:t> run in the document; two\n# French labels (\"Mat\u00e9riel\", \"Invitation \u00e0 la discussion\") occur twice and\n# both occurrences must be replaced, so wdReplaceAll is used for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Titre de la vid\u00e9o\", \"Titolo del video\"),\n    @(\"Sujet\", \"Argomento\"),\n    @(\"G\u00e9om\u00e9trie\", \"Geometria\"),\n    @(\"Objectif(s)\", \"Obiettivo/i\"),\n    @(\"Dur\u00e9e\", \"Lunghezza\"),\n    @(\"Localisation du camp\", \"Posizione del Campo\"),\n    @(\"Assistant\u00b7e\u00b7s\", \"Facilitatori\"),\n    @(\"N. d'\u00e9tudiant\u00b7e\u00b7s\", \"N. di studenti\"),\n    @(\"Date\", \"Data\"),\n    @(\"Mat\u00e9riel\", \"Materiale\"),\n    @(\"requis\", \"necessario\"),\n    @(\"Pr\u00e9paration\", \"Preparazione\"),\n    @(\"Temps vid\u00e9o\", \"Tempo del Video\"),\n    @(\"Introduction de la vid\u00e9o\", \"Introduzione al video\"),\n    @(\"Introduction de la premi\u00e8re activit\u00e9\", \"Introduzione all'attivit\u00e0 principale\"),\n    @(\"Les 8 graphes dessin\u00e9s pour la vid\u00e9o\", \"Gli 8 grafici disegnati per il video\"),\n    @(\"Introduction de l'exp\u00e9rience\", \"Introduzione all'esperimento\"),\n    @(\"Invitation \u00e0 la discussion\", \"Invito alla discussione\"),\n    @(\"Conclusion\", \"Conclusione\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
